$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$changes = @{
    "E4"   = 16.495
    "E6"   = 16.322
    "E7"   = 16.48
    "E16"  = 16.618
    "E20"  = 16.24
    "E28"  = 17.006
    "E29"  = 16.941
    "E32"  = 17.4
    "E40"  = 16.614
    "E46"  = 16.858
    "E51"  = 16.765
    "E52"  = 16.826
    "E57"  = 16.587
    "E59"  = 16.482
    "E62"  = 16.588
    "E66"  = 17.482
    "E73"  = 16.439
    "E74"  = 16.484
    "E92"  = 17.867
    "E100" = 16.768
}

foreach ($addr in $changes.Keys) {
    $ws.Range($addr).Value = $changes[$addr]
}
